$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values scraped for the cryptos list update (commit 2024-02-26).
# Cells whose new text would otherwise be auto-parsed as a number (e.g. "396.71")
# are forced to Text format first so they round-trip as strings, matching the source data.

$ws.Range('D2').Value = '53.387.72'
$ws.Range('E2').Value = '  +3.60%  '
$ws.Range('D3').Value = '3.136.69'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '396.71'
$ws.Range('E5').Value = '  +2.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.01'
$ws.Range('E6').Value = '  +5.80%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.545'
$ws.Range('E7').Value = '  +0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  +4.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.75'
$ws.Range('E10').Value = '  +5.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.139'
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0870'
$ws.Range('E12').Value = '  +1.16%  '
$ws.Range('D13').Value = '3.649.16'
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.03'
$ws.Range('E14').Value = '  +2.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.98'
$ws.Range('E15').Value = '  +2.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.05'
$ws.Range('E16').Value = '  +8.30%  '
$ws.Range('D17').Value = '3.149.43'
$ws.Range('E17').Value = '  +2.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.51'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('D19').Value = '53.384.05'
$ws.Range('E19').Value = '  +3.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.26'
$ws.Range('E20').Value = '  +3.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.73'
$ws.Range('E21').Value = '  +2.42%  '
$ws.Range('D22').Value = '0.0₃0972'
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.87'
$ws.Range('E23').Value = '  +0.99%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '270.65'
$ws.Range('E24').Value = '  +1.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.24'
$ws.Range('E25').Value = '  +2.99%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.97'
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.39'
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.27'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.169'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.110'
$ws.Range('E31').Value = '  +2.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '10.95'
$ws.Range('E32').Value = '  +6.62%  '
$ws.Range('B33').Value = 'VeChain'
$ws.Range('C33').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0499'
$ws.Range('E33').Value = '  +10.58%  '
$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '37.07'
$ws.Range('E34').Value = '  +6.65%  '
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '50.40'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.64'
$ws.Range('E37').Value = '  +9.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.78'
$ws.Range('E39').Value = '  +8.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.09'
$ws.Range('E40').Value = '  +9.14%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.289'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.19'
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '129.87'
$ws.Range('E44').Value = '  +3.65%  '
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.12'
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('D49').Value = '2.076.23'
$ws.Range('E49').Value = '  +2.12%  '
$ws.Range('B50').Value = 'BEAM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0337'
$ws.Range('E50').Value = '  +6.01%  '
$ws.Range('B51').Value = 'FlareNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0503'
$ws.Range('E51').Value = '  +16.84%  '
